# Rebuild DE_table1_dim10 / Sheet1 with the exp3 results: columns are
# reordered to F#, Best, Worst, Median, Mean, Std, Success Rate and the
# function rows are re-sorted alphabetically (F14, Fc1, Fc2, Fc6, Fc7, Fc9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "F#"
$ws.Cells.Item(1, 2).Value = "Best"
$ws.Cells.Item(1, 3).Value = "Worst"
$ws.Cells.Item(1, 4).Value = "Median"
$ws.Cells.Item(1, 5).Value = "Mean"
$ws.Cells.Item(1, 6).Value = "Std"
$ws.Cells.Item(1, 7).Value = "Success Rate"

# Data rows
$ws.Cells.Item(2, 1).Value = "F14"
$ws.Cells.Item(2, 2).Value = 0.067544
$ws.Cells.Item(2, 3).Value = 27.817729
$ws.Cells.Item(2, 4).Value = 1.298433
$ws.Cells.Item(2, 5).Value = 3.395156
$ws.Cells.Item(2, 6).Value = 4.742687
$ws.Cells.Item(2, 7).Value = 0

$ws.Cells.Item(3, 1).Value = "Fc1"
$ws.Cells.Item(3, 2).Value = 218822.679148
$ws.Cells.Item(3, 3).Value = 71395048.680576
$ws.Cells.Item(3, 4).Value = 9218706.618487
$ws.Cells.Item(3, 5).Value = 14620604.211877
$ws.Cells.Item(3, 6).Value = 16184793.225555
$ws.Cells.Item(3, 7).Value = 0

$ws.Cells.Item(4, 1).Value = "Fc2"
$ws.Cells.Item(4, 2).Value = 4501341.472907
$ws.Cells.Item(4, 3).Value = 2520542316.717356
$ws.Cells.Item(4, 4).Value = 275066199.127894
$ws.Cells.Item(4, 5).Value = 539937960.049946
$ws.Cells.Item(4, 6).Value = 627868863.143072
$ws.Cells.Item(4, 7).Value = 0

$ws.Cells.Item(5, 1).Value = "Fc6"
$ws.Cells.Item(5, 2).Value = 1.305127
$ws.Cells.Item(5, 3).Value = 8.077939
$ws.Cells.Item(5, 4).Value = 5.512057
$ws.Cells.Item(5, 5).Value = 5.302732
$ws.Cells.Item(5, 6).Value = 1.684117
$ws.Cells.Item(5, 7).Value = 0

$ws.Cells.Item(6, 1).Value = "Fc7"
$ws.Cells.Item(6, 2).Value = 0.913873
$ws.Cells.Item(6, 3).Value = 83.633985
$ws.Cells.Item(6, 4).Value = 15.133886
$ws.Cells.Item(6, 5).Value = 17.383763
$ws.Cells.Item(6, 6).Value = 15.360749
$ws.Cells.Item(6, 7).Value = 0

$ws.Cells.Item(7, 1).Value = "Fc9"
$ws.Cells.Item(7, 2).Value = 7.838607
$ws.Cells.Item(7, 3).Value = 60.823539
$ws.Cells.Item(7, 4).Value = 27.755326
$ws.Cells.Item(7, 5).Value = 28.986674
$ws.Cells.Item(7, 6).Value = 12.500881
$ws.Cells.Item(7, 7).Value = 0

